# changed MP time limit and corrected error in fixed recourse data
#
# Sheet1 (summary, rows 2-11): corrected "objective"/"solve time" values and
# updated the MP time-limit-driven counters (num cuts/num variables/num
# cons/num quad cons: F/G/H/I columns) from the 1/272/245/20 run to the
# 20/4490/4900/400 run.
#
# Sheets "1".."10" (per-instance CCG iteration logs): updated the second
# iteration's objective/MP gap/MP solve time (and the first iteration's
# MP solve time / Worst violation) to match the corrected fixed-recourse
# data produced after the MP time limit change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Cells.Item(2, 2).Value = -71.65924129199975
$ws.Cells.Item(2, 3).Value = 25.859401074
$ws.Cells.Item(2, 6).Value = 20
$ws.Cells.Item(2, 7).Value = 4490
$ws.Cells.Item(2, 8).Value = 4900
$ws.Cells.Item(2, 9).Value = 400
$ws.Cells.Item(3, 2).Value = -69.4525623469496
$ws.Cells.Item(3, 3).Value = 3.241466568
$ws.Cells.Item(3, 6).Value = 20
$ws.Cells.Item(3, 7).Value = 4490
$ws.Cells.Item(3, 8).Value = 4900
$ws.Cells.Item(3, 9).Value = 400
$ws.Cells.Item(4, 2).Value = -71.24685991078216
$ws.Cells.Item(4, 3).Value = 28.250229126
$ws.Cells.Item(4, 6).Value = 20
$ws.Cells.Item(4, 7).Value = 4490
$ws.Cells.Item(4, 8).Value = 4900
$ws.Cells.Item(4, 9).Value = 400
$ws.Cells.Item(5, 2).Value = -71.65440618645337
$ws.Cells.Item(5, 3).Value = 12.416208294
$ws.Cells.Item(5, 6).Value = 20
$ws.Cells.Item(5, 7).Value = 4490
$ws.Cells.Item(5, 8).Value = 4900
$ws.Cells.Item(5, 9).Value = 400
$ws.Cells.Item(6, 2).Value = -70.22823721822225
$ws.Cells.Item(6, 3).Value = 3.259004692
$ws.Cells.Item(6, 6).Value = 20
$ws.Cells.Item(6, 7).Value = 4490
$ws.Cells.Item(6, 8).Value = 4900
$ws.Cells.Item(6, 9).Value = 400
$ws.Cells.Item(7, 2).Value = -71.84370845876089
$ws.Cells.Item(7, 3).Value = 10.915439388
$ws.Cells.Item(7, 6).Value = 20
$ws.Cells.Item(7, 7).Value = 4490
$ws.Cells.Item(7, 8).Value = 4900
$ws.Cells.Item(7, 9).Value = 400
$ws.Cells.Item(8, 2).Value = -67.12860513254925
$ws.Cells.Item(8, 3).Value = 14.342945844
$ws.Cells.Item(8, 6).Value = 20
$ws.Cells.Item(8, 7).Value = 4490
$ws.Cells.Item(8, 8).Value = 4900
$ws.Cells.Item(8, 9).Value = 400
$ws.Cells.Item(9, 2).Value = -71.56787798834289
$ws.Cells.Item(9, 3).Value = 1.532547922
$ws.Cells.Item(9, 6).Value = 20
$ws.Cells.Item(9, 7).Value = 4490
$ws.Cells.Item(9, 8).Value = 4900
$ws.Cells.Item(9, 9).Value = 400
$ws.Cells.Item(10, 2).Value = -70.884938041055
$ws.Cells.Item(10, 3).Value = 4.897122523
$ws.Cells.Item(10, 6).Value = 20
$ws.Cells.Item(10, 7).Value = 4490
$ws.Cells.Item(10, 8).Value = 4900
$ws.Cells.Item(10, 9).Value = 400
$ws.Cells.Item(11, 2).Value = -68.19423108211265
$ws.Cells.Item(11, 3).Value = 1.227735495
$ws.Cells.Item(11, 6).Value = 20
$ws.Cells.Item(11, 7).Value = 4490
$ws.Cells.Item(11, 8).Value = 4900
$ws.Cells.Item(11, 9).Value = 400

$ws = $wb.Worksheets.Item("9")
$ws.Cells.Item(2, 4).Value = 0.031195961033447267
$ws.Cells.Item(2, 5).Value = 41.2321
$ws.Cells.Item(3, 2).Value = -70.884938041055
$ws.Cells.Item(3, 4).Value = 4.67280581173706

$ws = $wb.Worksheets.Item("10")
$ws.Cells.Item(2, 4).Value = 0.019166825015014648
$ws.Cells.Item(2, 5).Value = 39.2631
$ws.Cells.Item(3, 2).Value = -68.19423108211265
$ws.Cells.Item(3, 3).Value = 0.09386039621602227
$ws.Cells.Item(3, 4).Value = 1.0458227377525635

$ws = $wb.Worksheets.Item("1")
$ws.Cells.Item(2, 4).Value = 1.1170247866313476
$ws.Cells.Item(2, 5).Value = 38.78663
$ws.Cells.Item(3, 2).Value = -71.65924129199975
$ws.Cells.Item(3, 3).Value = 0.041979626255819115
$ws.Cells.Item(3, 4).Value = 13.638601026383789

$ws = $wb.Worksheets.Item("2")
$ws.Cells.Item(2, 4).Value = 0.012244432326538086
$ws.Cells.Item(2, 5).Value = 40.41718
$ws.Cells.Item(3, 2).Value = -69.4525623469496
$ws.Cells.Item(3, 3).Value = 0.03783807411897675
$ws.Cells.Item(3, 4).Value = 3.0640088524993896

$ws = $wb.Worksheets.Item("3")
$ws.Cells.Item(2, 4).Value = 0.014223903474975586
$ws.Cells.Item(2, 5).Value = 38.10258
$ws.Cells.Item(3, 2).Value = -71.24685991078216
$ws.Cells.Item(3, 3).Value = 0.05347504931540876
$ws.Cells.Item(3, 4).Value = 28.076176032256104

$ws = $wb.Worksheets.Item("4")
$ws.Cells.Item(2, 4).Value = 0.01347002674633789
$ws.Cells.Item(2, 5).Value = 40.53027
$ws.Cells.Item(3, 2).Value = -71.65440618645337
$ws.Cells.Item(3, 3).Value = 0.0011955544123725794
$ws.Cells.Item(3, 4).Value = 12.289500115434814

$ws = $wb.Worksheets.Item("5")
$ws.Cells.Item(2, 4).Value = 0.01768573743811035
$ws.Cells.Item(2, 5).Value = 41.85082
$ws.Cells.Item(3, 2).Value = -70.22823721822225
$ws.Cells.Item(3, 4).Value = 3.0841074175334473

$ws = $wb.Worksheets.Item("6")
$ws.Cells.Item(2, 4).Value = 0.022147351282714844
$ws.Cells.Item(2, 5).Value = 40.92989
$ws.Cells.Item(3, 2).Value = -71.84370845876089
$ws.Cells.Item(3, 4).Value = 10.738329696549316

$ws = $wb.Worksheets.Item("7")
$ws.Cells.Item(2, 4).Value = 0.027114733259887697
$ws.Cells.Item(2, 5).Value = 36.76141
$ws.Cells.Item(3, 2).Value = -67.12860513254925
$ws.Cells.Item(3, 3).Value = 0.06488405414626831
$ws.Cells.Item(3, 4).Value = 14.180363214760376

$ws = $wb.Worksheets.Item("8")
$ws.Cells.Item(2, 4).Value = 0.011153691938598633
$ws.Cells.Item(2, 5).Value = 41.44521
$ws.Cells.Item(3, 2).Value = -71.56787798834289
$ws.Cells.Item(3, 3).Value = 0.014902672757292004
$ws.Cells.Item(3, 4).Value = 1.3956988863045654
